$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 and Row 9 swap region names (West Midlands <-> Yorkshire and The Humber)
$ws.Range("A8").Value = "Yorkshire and The Humber"
$ws.Range("A9").Value = "West Midlands"

# Update Freq (column B) and PC (column C) values for rows 2-14
$ws.Range("B2").Value = 622
$ws.Range("C2").Value = 14.71

$ws.Range("B3").Value = 589
$ws.Range("C3").Value = 13.93

$ws.Range("B4").Value = 562
$ws.Range("C4").Value = 13.29

$ws.Range("B5").Value = 414
$ws.Range("C5").Value = 9.789999999999999

$ws.Range("B6").Value = 325
$ws.Range("C6").Value = 7.69

$ws.Range("B7").Value = 317
$ws.Range("C7").Value = 7.5

$ws.Range("B8").Value = 287
$ws.Range("C8").Value = 6.79

$ws.Range("B9").Value = 283
$ws.Range("C9").Value = 6.69

$ws.Range("B10").Value = 276
$ws.Range("C10").Value = 6.53

$ws.Range("B11").Value = 259
$ws.Range("C11").Value = 6.12

$ws.Range("B12").Value = 136
$ws.Range("C12").Value = 3.22

$ws.Range("B13").Value = 107
$ws.Range("C13").Value = 2.53

$ws.Range("B14").Value = 52
$ws.Range("C14").Value = 1.23
